$wb = $excel.ActiveWorkbook

# 1. Rename the existing "Sample" sheet to "list"
$listSheet = $wb.Worksheets.Item("Sample")
$listSheet.Name = "list"

# 2. Add a brand new worksheet right after "list", named "sample data"
$dataSheet = $wb.Worksheets.Add($null, $listSheet)
$dataSheet.Name = "sample data"

# 3. Populate the new sheet with form-like KEY/VAL data (order matters for
#    shared-string table layout: KEY/VAL header row, then labels down
#    column A, then values down column B)
$dataSheet.Range("A1").Value = "KEY"
$dataSheet.Range("B1").Value = "VAL"
$dataSheet.Range("A1:B1").Font.Bold = $true

$dataSheet.Range("A2").Value = "Email address"
$dataSheet.Range("A3").Value = "Password"

$dataSheet.Range("B2").Value = "root@localhost.localdomain"
$dataSheet.Range("B3").Value = "unsecure"

# Hyperlink on B2 (mailto link to the email address)
$dataSheet.Hyperlinks.Add($dataSheet.Range("B2"), "mailto:root@localhost.localdomain")

# Column width to fit the longest entry
$dataSheet.Columns.Item(1).ColumnWidth = 11.666666666666666

# Match the page margins used on the rest of the workbook
$dataSheet.PageSetup.LeftMargin = 54
$dataSheet.PageSetup.RightMargin = 54
$dataSheet.PageSetup.TopMargin = 72
$dataSheet.PageSetup.BottomMargin = 72
$dataSheet.PageSetup.HeaderMargin = 36
$dataSheet.PageSetup.FooterMargin = 36

# Select row 4 (below the data) then activate the sheet, matching the
# author's final selection/active-tab state
$dataSheet.Rows.Item(4).Select()
$dataSheet.Activate()
